# Split the single "Bibliografia" paragraph run into multiple runs
# separated by manual line breaks (<w:br/>), one before each
# lettered reference (B), C), D), E), F)).

$d = $word.ActiveDocument

function Insert-LineBreakBefore($anchorText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Text = $anchorText
    $find.Replacement.Text = "^l" + $anchorText
    $find.Forward = $true
    $find.Wrap = 0
    $find.Format = $false
    $find.MatchCase = $false
    $find.MatchWholeWord = $false
    $find.MatchWildcards = $false
    $find.MatchSoundsLike = $false
    $find.MatchAllWordForms = $false
    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)
}

Insert-LineBreakBefore "B)GRUS"
Insert-LineBreakBefore "C)VanderPlas"
Insert-LineBreakBefore "D)BANIN"
Insert-LineBreakBefore "E)BARI"
Insert-LineBreakBefore "F)Manuais"
